# Update countries & provincias Spain
# Applies the periodic COVID data refresh + three country-row swaps
# (Alemania/Brasil, Kenia/Mali, Curazao/Dominica) and bumps the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos (values refreshed, country unchanged) ---
$ws.Range("B4").Value = 1402275
$ws.Range("C4").Value = 16441
$ws.Range("D4").Value = 276175
$ws.Range("E4").Value = 1042979
$ws.Range("F4").Value = 16445
$ws.Range("G4").Value = 1326
$ws.Range("H4").Value = 83121

# --- Rows 10/11: Alemania and Brasil swap places, values refreshed ---
$ws.Range("A10").Value = "Brasil"
$ws.Range("B10").Value = 173126
$ws.Range("C10").Value = 3983
$ws.Range("D10").Value = 67384
$ws.Range("E10").Value = 93680
$ws.Range("F10").Value = 8318
$ws.Range("G10").Value = 437
$ws.Range("H10").Value = 12062

$ws.Range("A11").Value = "Alemania"
$ws.Range("B11").Value = 173031
$ws.Range("C11").Value = 455
$ws.Range("D11").Value = 147200
$ws.Range("E11").Value = 18114
$ws.Range("F11").Value = 1539
$ws.Range("G11").Value = 56
$ws.Range("H11").Value = 7717

# --- Row 15: India (values refreshed) ---
$ws.Range("D15").Value = 24420
$ws.Range("E15").Value = 47408

# --- Row 17: Canada (values refreshed) ---
$ws.Range("B17").Value = 71100
$ws.Range("C17").Value = 1119
$ws.Range("D17").Value = 33711
$ws.Range("E17").Value = 32222

# --- Row 52: Noruega (values refreshed) ---
$ws.Range("B52").Value = 8152
$ws.Range("C52").Value = 20
$ws.Range("E52").Value = 7892
$ws.Range("F52").Value = 23
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 228

# --- Row 100: El Salvador (values refreshed) ---
$ws.Range("F100").Value = 15

# --- Rows 114/115: Kenia and Mali swap places, values refreshed ---
$ws.Range("A114").Value = "Mali"
$ws.Range("B114").Value = 730
$ws.Range("C114").Value = 18
$ws.Range("D114").Value = 398
$ws.Range("E114").Value = 292
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 40

$ws.Range("A115").Value = "Kenia"
$ws.Range("B115").Value = 715
$ws.Range("C115").Value = 15
$ws.Range("D115").Value = 259
$ws.Range("E115").Value = 420
$ws.Range("F115").Value = 1
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 36

# --- Row 151: Birmania (values refreshed) ---
$ws.Range("D151").Value = 76
$ws.Range("E151").Value = 98

# --- Row 158: Uganda (values refreshed) ---
$ws.Range("B158").Value = 129
$ws.Range("C158").Value = 8
$ws.Range("E158").Value = 74

# --- Rows 198/199: Curazao and Dominica swap places ---
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

# --- Timestamp string update ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 22:05"
